# Add a "Profession" column to the Data sheet, right after "Caregiver"
# and before "Specialty", with value "RN" for every data row.
#
# Current layout (before insert):
# A Facility | B Cost Center/Type | C Unit | D Subcontractor | E Caregiver
# F Specialty | G Date | H Shift | I Task Name | J Regular Hours | K OT Hours
# L Holiday Hours | M Total Hours | N Bill Rate | O Bill Gross | P Invoice #
# Q Invoicing Period | R Total Hrs/Caregiver | S Total Bill Gross/Caregiver
#
# Insert a new column F so "Profession" lands between Caregiver (E) and
# Specialty (which becomes G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("F1").EntireColumn.Insert()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Cells.Item(1, 6).Value = "Profession"
$ws.Cells.Item(1, 6).Style = $ws.Cells.Item(1, 5).Style

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "RN"
    $ws.Cells.Item($r, 6).Style = $ws.Cells.Item($r, 5).Style
}
